# MenuInterativo.xlsx — add a "Designer Partner" entry as the new first
# data row (row 2) of the menu table, pushing the existing entries
# (Hardwares, Instalando programas automaticamente!, ...) down by one row.
# This also shifts the dependent "id=" helper table (rows 16-29) down to
# rows 17-30, since it lives on the same sheet below the main table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (Hardwares) into a fresh row above it so the new row
# inherits all of row 2's number formats / styles (including the blank
# G2 style used by the "</i>" spacer column), then the old content that
# was in row 2 slides down to row 3, etc.
$ws.Rows("2").Copy()
$ws.Rows("2").Insert()

# The insert logic picks up the style of the row above (row 1) for the
# blank G column; restore it to match the rest of the table (same blank
# style as G3).
$ws.Range("G3").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# Fill in the new row's content: "Designer Partner".
$ws.Range("A2").Value = "Designer Partner"
$ws.Range("B2").Formula = "=A2"
$ws.Range("D2").Formula = "=B2"
$ws.Range("H2").Formula = "=A2"

# Match the author's final selection.
$ws.Range("D2").Select()
